$d = $word.ActiveDocument

$pairs = @(
    @("96×82=", "17×47="),
    @("32×74=", "12×60="),
    @("85×18=", "27×45="),
    @("20×71=", "29×23="),
    @("93×13=", "95×12="),
    @("96×67=", "79×12="),
    @("88×41=", "25×23="),
    @("71×61=", "42×63="),
    @("60×98=", "50×73="),
    @("63×34=", "35×28="),
    @("11×28=", "99×92="),
    @("86×25=", "89×99="),
    @("62×89=", "78×51="),
    @("31×77=", "74×88="),
    @("48×63=", "37×19="),
    @("44×23=", "20×40="),
    @("18×91=", "45×92="),
    @("36×63=", "61×97="),
    @("58×46=", "11×60="),
    @("15×24=", "66×85="),
    @("92×56=", "26×47="),
    @("15×99=", "60×72="),
    @("54×41=", "79×67="),
    @("90×18=", "33×58="),
    @("29×11=", "42×95=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
